$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 586.5714
$ws.Range("I4").Value = 624
$ws.Range("K4").Value = 624
$ws.Range("M4").Value = -510
$ws.Range("H41").Value = 1474.75
$ws.Range("I41").Value = 1466.6666
$ws.Range("J41").Value = 1499
$ws.Range("K41").Value = 1466.6666
$ws.Range("L41").Value = 1499
$ws.Range("M41").Value = -1026.6666
$ws.Range("N41").Value = -2379
$ws.Range("H51").Value = 17533.7
$ws.Range("I51").Value = 109990
$ws.Range("J51").Value = 7260.778
$ws.Range("K51").Value = 109990
$ws.Range("L51").Value = 7260.778
$ws.Range("M51").Value = -109506
$ws.Range("N51").Value = -8228.778
$ws.Range("H62").Value = 7942301.5
$ws.Range("I62").Value = 12994511
$ws.Range("K62").Value = 12994511
$ws.Range("M62").Value = -12993887
$ws.Range("H65").Value = 7942301.5
$ws.Range("I65").Value = 12994511
$ws.Range("K65").Value = 64972555
$ws.Range("M65").Value = -64969435
$ws.Range("H70").Value = 1333.7142
$ws.Range("I70").Value = 2349
$ws.Range("J70").Value = 1164.5
$ws.Range("K70").Value = 7047
$ws.Range("L70").Value = 3493.5
$ws.Range("M70").Value = -6777
$ws.Range("N70").Value = -4033.5
$ws.Range("H73").Value = 1333.7142
$ws.Range("I73").Value = 2349
$ws.Range("J73").Value = 1164.5
$ws.Range("K73").Value = 7047
$ws.Range("L73").Value = 3493.5
$ws.Range("M73").Value = -6111
$ws.Range("N73").Value = -5365.5
$ws.Range("H76").Value = 62506676
$ws.Range("I76").Value = 100005816
$ws.Range("K76").Value = 100005816
$ws.Range("M76").Value = -100005501
$ws.Range("H79").Value = 62506676
$ws.Range("I79").Value = 100005816
$ws.Range("K79").Value = 100005816
$ws.Range("M79").Value = -100004724
$ws.Range("H86").Value = 2938.125
$ws.Range("I86").Value = 3478.25
$ws.Range("J86").Value = 2398
$ws.Range("K86").Value = 3478.25
$ws.Range("L86").Value = 2398
$ws.Range("M86").Value = -2355.25
$ws.Range("N86").Value = -4644
$ws.Range("H88").Value = 3173
$ws.Range("J88").Value = 3097.9092
$ws.Range("L88").Value = 3097.9092
$ws.Range("N88").Value = -3909.9092
$ws.Range("H89").Value = 2938.125
$ws.Range("I89").Value = 3478.25
$ws.Range("J89").Value = 2398
$ws.Range("K89").Value = 17391.25
$ws.Range("L89").Value = 11990
$ws.Range("M89").Value = -11775.25
$ws.Range("N89").Value = -23222
$ws.Range("H91").Value = 3173
$ws.Range("J91").Value = 3097.9092
$ws.Range("L91").Value = 3097.9092
$ws.Range("N91").Value = -5905.9092
$ws.Range("H98").Value = 1379.8518
$ws.Range("I98").Value = 1479.3334
$ws.Range("J98").Value = 1031.6666
$ws.Range("K98").Value = 1479.3334
$ws.Range("L98").Value = 1031.6666
$ws.Range("M98").Value = 18.66660000000002
$ws.Range("N98").Value = -4027.6666
$ws.Range("H106").Value = 4498.3
$ws.Range("I106").Value = 4453.722
$ws.Range("K106").Value = 4453.722
$ws.Range("M106").Value = -3822.722
$ws.Range("H107").Value = 3634.9285
$ws.Range("I107").Value = 2077.7778
$ws.Range("J107").Value = 6437.8
$ws.Range("K107").Value = 2077.7778
$ws.Range("L107").Value = 6437.8
$ws.Range("M107").Value = -157.7777999999998
$ws.Range("N107").Value = -10277.8
$ws.Range("H112").Value = 3581.6924
$ws.Range("J112").Value = 2296.8333
$ws.Range("L112").Value = 6890.499899999999
$ws.Range("N112").Value = -9106.499899999999
$ws.Range("H122").Value = 1379.8518
$ws.Range("I122").Value = 1479.3334
$ws.Range("J122").Value = 1031.6666
$ws.Range("K122").Value = 4438.0002
$ws.Range("L122").Value = 3094.9998
$ws.Range("M122").Value = -1988.0002
$ws.Range("N122").Value = -7994.9998
$ws.Range("H132").Value = 498541.06
$ws.Range("I132").Value = 646058.4399999999
$ws.Range("K132").Value = 1938175.32
$ws.Range("M132").Value = -1935645.32
$ws.Range("H137").Value = 2899.6667
$ws.Range("I137").Value = 2105.88
$ws.Range("J137").Value = 4703.727
$ws.Range("K137").Value = 6317.64
$ws.Range("L137").Value = 14111.181
$ws.Range("M137").Value = -3767.64
$ws.Range("N137").Value = -19211.181
$ws.Range("H141").Value = 2763.92
$ws.Range("I141").Value = 2795.75
$ws.Range("K141").Value = 8387.25
$ws.Range("M141").Value = -3207.25

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4633866.5
$ws.Range("I32").Value = 1932.25
$ws.Range("K32").Value = 1932.25
$ws.Range("M32").Value = -1645.25
$ws.Range("H74").Value = 4416.25
$ws.Range("I74").Value = 2224.75
$ws.Range("J74").Value = 5512
$ws.Range("K74").Value = 2224.75
$ws.Range("L74").Value = 5512
$ws.Range("M74").Value = -1350.75
$ws.Range("N74").Value = -7260
$ws.Range("H77").Value = 4416.25
$ws.Range("I77").Value = 2224.75
$ws.Range("J77").Value = 5512
$ws.Range("K77").Value = 11123.75
$ws.Range("L77").Value = 27560
$ws.Range("M77").Value = -6755.75
$ws.Range("N77").Value = -36296
$ws.Range("H102").Value = 1768.6111
$ws.Range("I102").Value = 1737.3529
$ws.Range("J102").Value = 2300
$ws.Range("K102").Value = 1737.3529
$ws.Range("L102").Value = 2300
$ws.Range("M102").Value = -115.3529000000001
$ws.Range("N102").Value = -5544
$ws.Range("H107").Value = 55000
$ws.Range("J107").Value = 55000
$ws.Range("L107").Value = 55000
$ws.Range("N107").Value = -62680

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 9679.25
$ws.Range("I24").Value = 9449.5
$ws.Range("K24").Value = 9449.5
$ws.Range("M24").Value = -9214.5
$ws.Range("H64").Value = 5103.3
$ws.Range("I64").Value = 1700
$ws.Range("J64").Value = 5954.125
$ws.Range("K64").Value = 1700
$ws.Range("L64").Value = 5954.125
$ws.Range("M64").Value = -1475
$ws.Range("N64").Value = -6404.125
$ws.Range("H67").Value = 5103.3
$ws.Range("I67").Value = 1700
$ws.Range("J67").Value = 5954.125
$ws.Range("K67").Value = 1700
$ws.Range("L67").Value = 5954.125
$ws.Range("M67").Value = -920
$ws.Range("N67").Value = -7514.125
$ws.Range("H96").Value = 25428
$ws.Range("I96").Value = 25428
$ws.Range("K96").Value = 25428
$ws.Range("M96").Value = -22682
$ws.Range("H134").Value = 1593432.9
$ws.Range("I134").Value = 1989323.9
$ws.Range("J134").Value = 9868.833000000001
$ws.Range("K134").Value = 5967971.699999999
$ws.Range("L134").Value = 29606.499
$ws.Range("M134").Value = -5965436.699999999
$ws.Range("N134").Value = -34676.499

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2387.3
$ws.Range("I31").Value = 2263.6667
$ws.Range("J31").Value = 3500
$ws.Range("K31").Value = 2263.6667
$ws.Range("L31").Value = 3500
$ws.Range("M31").Value = -1968.6667
$ws.Range("N31").Value = -4090
$ws.Range("H34").Value = 2387.3
$ws.Range("I34").Value = 2263.6667
$ws.Range("J34").Value = 3500
$ws.Range("K34").Value = 2263.6667
$ws.Range("L34").Value = 3500
$ws.Range("M34").Value = -2061.6667
$ws.Range("N34").Value = -3904
$ws.Range("H58").Value = 58827652
$ws.Range("I58").Value = 100001860
$ws.Range("J58").Value = 7358.7144
$ws.Range("K58").Value = 100001860
$ws.Range("L58").Value = 7358.7144
$ws.Range("M58").Value = -100001657
$ws.Range("N58").Value = -7764.7144
$ws.Range("H86").Value = 17076.666
$ws.Range("I86").Value = 27610
$ws.Range("J86").Value = 11810
$ws.Range("K86").Value = 27610
$ws.Range("L86").Value = 11810
$ws.Range("M86").Value = -26487
$ws.Range("N86").Value = -14056
$ws.Range("H89").Value = 17076.666
$ws.Range("I89").Value = 27610
$ws.Range("J89").Value = 11810
$ws.Range("K89").Value = 138050
$ws.Range("L89").Value = 59050
$ws.Range("M89").Value = -132434
$ws.Range("N89").Value = -70282
$ws.Range("H134").Value = 32267586
$ws.Range("I134").Value = 142861620
$ws.Range("J134").Value = 10993.792
$ws.Range("K134").Value = 428584860
$ws.Range("L134").Value = 32981.376
$ws.Range("M134").Value = -428582325
$ws.Range("N134").Value = -38051.376
$ws.Range("H136").Value = 58827652
$ws.Range("I136").Value = 100001860
$ws.Range("J136").Value = 7358.7144
$ws.Range("K136").Value = 300005580
$ws.Range("L136").Value = 22076.1432
$ws.Range("M136").Value = -300003030
$ws.Range("N136").Value = -27176.1432

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 219996.5
$ws.Range("J37").Value = 219996.5
$ws.Range("L37").Value = 659989.5
$ws.Range("N37").Value = -660213.5
$ws.Range("H63").Value = 12717
$ws.Range("I63").Value = 7108.1113
$ws.Range("K63").Value = 21324.3339
$ws.Range("M63").Value = -20575.3339
$ws.Range("H66").Value = 12717
$ws.Range("I66").Value = 7108.1113
$ws.Range("K66").Value = 63973.00169999999
$ws.Range("M66").Value = -60229.00169999999
$ws.Range("H87").Value = 18104.4
$ws.Range("I87").Value = 10507.333
$ws.Range("J87").Value = 29500
$ws.Range("K87").Value = 31521.999
$ws.Range("L87").Value = 88500
$ws.Range("M87").Value = -30273.999
$ws.Range("N87").Value = -90996
$ws.Range("H90").Value = 18104.4
$ws.Range("I90").Value = 10507.333
$ws.Range("J90").Value = 29500
$ws.Range("K90").Value = 94565.997
$ws.Range("L90").Value = 265500
$ws.Range("M90").Value = -88325.997
$ws.Range("N90").Value = -277980
$ws.Range("H98").Value = 2432.8333
$ws.Range("J98").Value = 1919.4
$ws.Range("L98").Value = 5758.200000000001
$ws.Range("N98").Value = -8754.200000000001
$ws.Range("H107").Value = 6702.45
$ws.Range("J107").Value = 7420.8887
$ws.Range("L107").Value = 22262.6661
$ws.Range("N107").Value = -26102.6661
$ws.Range("H117").Value = 1878.2307
$ws.Range("I117").Value = 3606.2
$ws.Range("K117").Value = 10818.6
$ws.Range("M117").Value = -7376.599999999999
$ws.Range("H129").Value = 11906527
$ws.Range("J129").Value = 20836120
$ws.Range("L129").Value = 62508360
$ws.Range("N129").Value = -62518360
$ws.Range("H140").Value = 43862430
$ws.Range("I140").Value = 64104324
$ws.Range("J140").Value = 4998.8335
$ws.Range("K140").Value = 192312972
$ws.Range("L140").Value = 14996.5005
$ws.Range("M140").Value = -192307792
$ws.Range("N140").Value = -25356.5005

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7058.8
$ws.Range("J113").Value = 9059.1
$ws.Range("L113").Value = 9059.1
$ws.Range("N113").Value = -13399.1
$ws.Range("H122").Value = 8533.299999999999
$ws.Range("I122").Value = 5832.3335
$ws.Range("J122").Value = 12584.75
$ws.Range("K122").Value = 17497.0005
$ws.Range("L122").Value = 37754.25
$ws.Range("M122").Value = -15047.0005
$ws.Range("N122").Value = -42654.25
$ws.Range("H126").Value = 35300120
$ws.Range("I126").Value = 100002940
$ws.Range("K126").Value = 300008820
$ws.Range("M126").Value = -300006350

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5167.5654
$ws.Range("I61").Value = 2778.6924
$ws.Range("K61").Value = 2778.6924
$ws.Range("M61").Value = -2576.6924
$ws.Range("H68").Value = 1925
$ws.Range("I68").Value = 1887.5
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1887.5
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1138.5
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1925
$ws.Range("I71").Value = 1887.5
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 9437.5
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -5693.5
$ws.Range("N71").Value = -17488
$ws.Range("H76").Value = 13737.5
$ws.Range("J76").Value = 12500
$ws.Range("L76").Value = 12500
$ws.Range("N76").Value = -13176
$ws.Range("H79").Value = 13737.5
$ws.Range("J79").Value = 12500
$ws.Range("L79").Value = 12500
$ws.Range("N79").Value = -14840
$ws.Range("H82").Value = 2555.2307
$ws.Range("I82").Value = 1027.4
$ws.Range("J82").Value = 4638.636
$ws.Range("K82").Value = 1027.4
$ws.Range("L82").Value = 4638.636
$ws.Range("M82").Value = -666.4000000000001
$ws.Range("N82").Value = -5360.636
$ws.Range("H85").Value = 2555.2307
$ws.Range("I85").Value = 1027.4
$ws.Range("J85").Value = 4638.636
$ws.Range("K85").Value = 1027.4
$ws.Range("L85").Value = 4638.636
$ws.Range("M85").Value = 220.5999999999999
$ws.Range("N85").Value = -7134.636
$ws.Range("H113").Value = 5167.5654
$ws.Range("I113").Value = 2778.6924
$ws.Range("K113").Value = 2778.6924
$ws.Range("M113").Value = -608.6923999999999
$ws.Range("H122").Value = 4837.0303
$ws.Range("I122").Value = 3209.4285
$ws.Range("K122").Value = 9628.2855
$ws.Range("M122").Value = -7178.2855
$ws.Range("H132").Value = 4247.375
$ws.Range("I132").Value = 4499.5
$ws.Range("K132").Value = 13498.5
$ws.Range("M132").Value = -10968.5
$ws.Range("H136").Value = 12199686
$ws.Range("I136").Value = 71429440
$ws.Range("K136").Value = 214288320
$ws.Range("M136").Value = -214285770

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 29682.334
$ws.Range("J52").Value = 23523.5
$ws.Range("L52").Value = 23523.5
$ws.Range("N52").Value = -23975.5
$ws.Range("H58").Value = 32663.334
$ws.Range("I58").Value = 32663.334
$ws.Range("K58").Value = 32663.334
$ws.Range("M58").Value = -32355.334
$ws.Range("H81").Value = 1494.091
$ws.Range("I81").Value = 1413.7
$ws.Range("K81").Value = 2827.4
$ws.Range("M81").Value = -1766.4
$ws.Range("H84").Value = 1494.091
$ws.Range("I84").Value = 1413.7
$ws.Range("K84").Value = 14137
$ws.Range("M84").Value = -8833
$ws.Range("H113").Value = 20635992
$ws.Range("I113").Value = 44219376
$ws.Range("J113").Value = 532.5
$ws.Range("K113").Value = 132658128
$ws.Range("L113").Value = 1597.5
$ws.Range("M113").Value = -132655958
$ws.Range("N113").Value = -5937.5
$ws.Range("H122").Value = 2801.3704
$ws.Range("I122").Value = 2365.8125
$ws.Range("K122").Value = 7097.4375
$ws.Range("M122").Value = -4647.4375
$ws.Range("H132").Value = 9803.286
$ws.Range("I132").Value = 4710
$ws.Range("J132").Value = 11840.6
$ws.Range("K132").Value = 14130
$ws.Range("L132").Value = 35521.8
$ws.Range("M132").Value = -11600
$ws.Range("N132").Value = -40581.8
